# ---------------------------------------------------------------------------
# Update docs/epexspot_prices.xlsx with the latest scraped EPEX spot / Gaz /
# CO2 prices:
#   * "Prix Spot" sheet: a new daily price column ("02-dec") is inserted
#     right before the existing "01-oct." column, shifting every following
#     date column one to the right (DY:FC -> DZ:FD). The new column gets
#     "-" placeholders for every hour row (no data yet for that date).
#   * "Gaz" sheet: two new rows are appended for 2025-11-29 / 2025-11-30.
#   * "CO2" sheet: two new rows are appended for 2025-11-29 / 2025-11-30
#     (price not published yet, so the price cell stays blank).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# xlPasteSpecial constants used below.
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1) "Prix Spot" sheet — insert a new date column before DY ("01-oct.")
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Column 129 == "DY". Inserting here pushes DY:FC (the old content) one
# column to the right, to DZ:FD, and leaves a blank DY column behind,
# inheriting the surrounding formatting (bold / centered / bordered header
# style for row 1, plain style for the data rows).
$ws1.Cells.Item(1, 129).EntireColumn.Insert()

# Header label for the freshly inserted column.
$ws1.Cells.Item(1, 129).Value = "02-dec"

# No observations yet for that date -> "-" placeholder for every hour row.
for ($r = 2; $r -le 25; $r++) {
    $ws1.Cells.Item($r, 129).Value = "-"
}

# ---------------------------------------------------------------------------
# 2) "Gaz" sheet — append 2025-11-29 and 2025-11-30
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Gaz")

# Column A holds dates stored as plain text (e.g. "2025-11-17"), not real
# Excel dates. Assigning a "YYYY-MM-DD" string straight to .Value gets
# auto-recognised as a date by Excel, so force Text formatting first, then
# copy the (already-Text) format from the previous row back on top so the
# cell ends up with no extra explicit style, matching its neighbours.
$gazPrev = $ws2.Cells.Item(155, 1)

$gazDate1 = $ws2.Cells.Item(156, 1)
$gazDate1.NumberFormat = "@"
$gazDate1.Value = "2025-11-29"
$gazPrev.Copy()
$gazDate1.PasteSpecial($xlPasteFormats)

$gazDate2 = $ws2.Cells.Item(157, 1)
$gazDate2.NumberFormat = "@"
$gazDate2.Value = "2025-11-30"
$gazPrev.Copy()
$gazDate2.PasteSpecial($xlPasteFormats)

$ws2.Cells.Item(156, 2).Value = 27.525
$ws2.Cells.Item(157, 2).Value = 27.525

# ---------------------------------------------------------------------------
# 3) "CO2" sheet — append 2025-11-29 and 2025-11-30 (price not out yet)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CO2")

$co2Prev = $ws3.Cells.Item(156, 1)

$co2Date1 = $ws3.Cells.Item(157, 1)
$co2Date1.NumberFormat = "@"
$co2Date1.Value = "2025-11-29"
$co2Prev.Copy()
$co2Date1.PasteSpecial($xlPasteFormats)

$co2Date2 = $ws3.Cells.Item(158, 1)
$co2Date2.NumberFormat = "@"
$co2Date2.Value = "2025-11-30"
$co2Prev.Copy()
$co2Date2.PasteSpecial($xlPasteFormats)

# Price column intentionally left blank for both new CO2 rows.
